# Auto-generated script to apply odds updates for rows 7, 9, 13
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G7").Value = 2.57
$ws.Range("H7").Value = 3.2
$ws.Range("I7").Value = 2.6
$ws.Range("J7").Value = 3.1
$ws.Range("K7").Value = 2.1
$ws.Range("L7").Value = 3.1
$ws.Range("R7").Value = 1.91
$ws.Range("U7").Value = 1.55
$ws.Range("V7").Value = 2.15
$ws.Range("W7").Value = 9.75
$ws.Range("X7").Value = 14.5
$ws.Range("Y7").Value = 9.5
$ws.Range("Z7").Value = 30
$ws.Range("AA7").Value = 20
$ws.Range("AB7").Value = 25
$ws.Range("AE7").Value = 11.5
$ws.Range("AG7").Value = 300
$ws.Range("AI7").Value = 15
$ws.Range("AJ7").Value = 9.5
$ws.Range("AK7").Value = 32
$ws.Range("AL7").Value = 20
$ws.Range("AM7").Value = 24
$ws.Range("AN7").Value = 4.6
$ws.Range("AO7").Value = 13.5
$ws.Range("AP7").Value = 19
$ws.Range("AQ7").Value = 55
$ws.Range("AR7").Value = 80
$ws.Range("AS7").Value = 200
$ws.Range("AW7").Value = 4.65
$ws.Range("AX7").Value = 13.5
$ws.Range("AY7").Value = 18.5
$ws.Range("AZ7").Value = 55
$ws.Range("BA7").Value = 80
$ws.Range("BB7").Value = 200
$ws.Range("J9").Value = 1.92
$ws.Range("G13").Value = 1.42
$ws.Range("H13").Value = 4.1
$ws.Range("I13").Value = 8.5
$ws.Range("M13").Value = 1.06
$ws.Range("N13").Value = 10
$ws.Range("W13").Value = 6.5
$ws.Range("X13").Value = 6.5
$ws.Range("Y13").Value = 9
$ws.Range("Z13").Value = 9
$ws.Range("AD13").Value = 8
$ws.Range("AH13").Value = 19
$ws.Range("AJ13").Value = 23
$ws.Range("AL13").Value = 51
$ws.Range("AM13").Value = 51
$ws.Range("AO13").Value = 7
$ws.Range("AQ13").Value = 21
$ws.Range("AU13").Value = 9.5
$ws.Range("AW13").Value = 8.5
$ws.Range("AZ13").Value = 151
